# repull data, push all data, mean calculation
# Update column F (dSF) values for rows 2-27 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -1
$ws.Range("F3").Value  = -3
$ws.Range("F4").Value  = -3
$ws.Range("F5").Value  = -4
$ws.Range("F6").Value  = -1
$ws.Range("F7").Value  = -3
$ws.Range("F8").Value  = -7
$ws.Range("F9").Value  = -1
$ws.Range("F10").Value = -1
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = -1
$ws.Range("F17").Value = -2
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = 6
$ws.Range("F22").Value = 2
$ws.Range("F24").Value = 4
$ws.Range("F25").Value = 4
$ws.Range("F26").Value = -3
$ws.Range("F27").Value = -3
